$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 346.46155  # H15: was 490.42856
$ws.Cells.Item(15, 9).Value = 346.46155  # I15: was 490.42856
$ws.Cells.Item(15, 11).Value = 1039.38465  # K15: was 1471.28568
$ws.Cells.Item(15, 13).Value = -870.38465  # M15: was -1302.28568
# Row 33
$ws.Cells.Item(33, 8).Value = 151.64285  # H33: was 164.92857
$ws.Cells.Item(33, 9).Value = 157.18182  # I33: was 162.23077
$ws.Cells.Item(33, 10).Value = 131.33333  # J33: was 200
$ws.Cells.Item(33, 11).Value = 157.18182  # K33: was 162.23077
$ws.Cells.Item(33, 12).Value = 131.33333  # L33: was 200
$ws.Cells.Item(33, 13).Value = 71.81818000000001  # M33: was 66.76922999999999
$ws.Cells.Item(33, 14).Value = -589.3333299999999  # N33: was -658
# Row 40
$ws.Cells.Item(40, 8).Value = 1906.762  # H40: was 1919.122
$ws.Cells.Item(40, 9).Value = 1702.9354  # I40: was 1713.0333
$ws.Cells.Item(40, 11).Value = 1702.9354  # K40: was 1713.0333
$ws.Cells.Item(40, 13).Value = -1527.9354  # M40: was -1538.0333
# Row 43
$ws.Cells.Item(43, 8).Value = 3333  # H43: was 0
$ws.Cells.Item(43, 10).Value = 3333  # J43: was 0
$ws.Cells.Item(43, 12).Value = 3333  # L43: was 0
$ws.Cells.Item(43, 14).Value = -3471  # N43: was None
# Row 51
$ws.Cells.Item(51, 8).Value = 4974.5  # H51: was 5000
$ws.Cells.Item(51, 9).Value = 4974.5  # I51: was 5000
$ws.Cells.Item(51, 11).Value = 4974.5  # K51: was 5000
$ws.Cells.Item(51, 13).Value = -4490.5  # M51: was -4516
# Row 58
$ws.Cells.Item(58, 8).Value = 822.2857  # H58: was 1166.6666
$ws.Cells.Item(58, 9).Value = 650  # I58: was 500
$ws.Cells.Item(58, 10).Value = 891.2  # J58: was 1500
$ws.Cells.Item(58, 11).Value = 1950  # K58: was 1500
$ws.Cells.Item(58, 12).Value = 2673.6  # L58: was 4500
$ws.Cells.Item(58, 13).Value = -1800  # M58: was -1350
$ws.Cells.Item(58, 14).Value = -2973.6  # N58: was -4800
# Row 64
$ws.Cells.Item(64, 8).Value = 4549.5  # H64: was 4832.1665
$ws.Cells.Item(64, 9).Value = 3649  # I64: was 3499
$ws.Cells.Item(64, 10).Value = 4999.75  # J64: was 5098.8
$ws.Cells.Item(64, 11).Value = 3649  # K64: was 3499
$ws.Cells.Item(64, 12).Value = 4999.75  # L64: was 5098.8
$ws.Cells.Item(64, 13).Value = -3401  # M64: was -3251
$ws.Cells.Item(64, 14).Value = -5495.75  # N64: was -5594.8
# Row 67
$ws.Cells.Item(67, 8).Value = 4549.5  # H67: was 4832.1665
$ws.Cells.Item(67, 9).Value = 3649  # I67: was 3499
$ws.Cells.Item(67, 10).Value = 4999.75  # J67: was 5098.8
$ws.Cells.Item(67, 11).Value = 3649  # K67: was 3499
$ws.Cells.Item(67, 12).Value = 4999.75  # L67: was 5098.8
$ws.Cells.Item(67, 13).Value = -2791  # M67: was -2641
$ws.Cells.Item(67, 14).Value = -6715.75  # N67: was -6814.8
# Row 88
$ws.Cells.Item(88, 8).Value = 1679.8182  # H88: was 1541.6923
$ws.Cells.Item(88, 9).Value = 1807.6  # I88: was 2084.75
$ws.Cells.Item(88, 10).Value = 1573.3334  # J88: was 1300.3334
$ws.Cells.Item(88, 11).Value = 1807.6  # K88: was 2084.75
$ws.Cells.Item(88, 12).Value = 1573.3334  # L88: was 1300.3334
$ws.Cells.Item(88, 13).Value = -1401.6  # M88: was -1678.75
$ws.Cells.Item(88, 14).Value = -2385.3334  # N88: was -2112.3334
# Row 91
$ws.Cells.Item(91, 8).Value = 1679.8182  # H91: was 1541.6923
$ws.Cells.Item(91, 9).Value = 1807.6  # I91: was 2084.75
$ws.Cells.Item(91, 10).Value = 1573.3334  # J91: was 1300.3334
$ws.Cells.Item(91, 11).Value = 1807.6  # K91: was 2084.75
$ws.Cells.Item(91, 12).Value = 1573.3334  # L91: was 1300.3334
$ws.Cells.Item(91, 13).Value = -403.5999999999999  # M91: was -680.75
$ws.Cells.Item(91, 14).Value = -4381.3334  # N91: was -4108.3334
# Row 107
$ws.Cells.Item(107, 8).Value = 0  # H107: was 415.66666
$ws.Cells.Item(107, 9).Value = 0  # I107: was 403.33334
$ws.Cells.Item(107, 10).Value = 0  # J107: was 428
$ws.Cells.Item(107, 11).Value = 0  # K107: was 403.33334
$ws.Cells.Item(107, 12).Value = 0  # L107: was 428
$ws.Cells.Item(107, 13).ClearContents()  # M107: was 1516.66666
$ws.Cells.Item(107, 14).ClearContents()  # N107: was -4268
# Row 132
$ws.Cells.Item(132, 8).Value = 1989.8235  # H132: was 2071.6875
$ws.Cells.Item(132, 9).Value = 1934.2142  # I132: was 2133.25
$ws.Cells.Item(132, 10).Value = 2249.3333  # J132: was 1887
$ws.Cells.Item(132, 11).Value = 5802.642599999999  # K132: was 6399.75
$ws.Cells.Item(132, 12).Value = 6747.999899999999  # L132: was 5661
$ws.Cells.Item(132, 13).Value = -3272.642599999999  # M132: was -3869.75
$ws.Cells.Item(132, 14).Value = -11807.9999  # N132: was -10721

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Cells.Item(97, 8).Value = 944.7273  # H97: was 889.7
$ws.Cells.Item(97, 9).Value = 890.7  # I97: was 823.55554
$ws.Cells.Item(97, 11).Value = 890.7  # K97: was 823.55554
$ws.Cells.Item(97, 13).Value = -394.7  # M97: was -327.55554
# Row 110
$ws.Cells.Item(110, 8).Value = 2847068.8  # H110: was 2643965.2
$ws.Cells.Item(110, 9).Value = 2847068.8  # I110: was 3364538.5
$ws.Cells.Item(110, 10).Value = 0  # J110: was 1862.6666
$ws.Cells.Item(110, 11).Value = 2847068.8  # K110: was 3364538.5
$ws.Cells.Item(110, 12).Value = 0  # L110: was 1862.6666
$ws.Cells.Item(110, 13).Value = -2845023.8  # M110: was -3362493.5
$ws.Cells.Item(110, 14).ClearContents()  # N110: was -5952.6666
# Row 122
$ws.Cells.Item(122, 8).Value = 1372.3334  # H122: was 1334.8
$ws.Cells.Item(122, 9).Value = 975.1667  # I122: was 978.2857
$ws.Cells.Item(122, 11).Value = 2925.5001  # K122: was 2934.8571
$ws.Cells.Item(122, 13).Value = -475.5001000000002  # M122: was -484.8571000000002
# Row 132
$ws.Cells.Item(132, 8).Value = 1270.3334  # H132: was 1270.6666
$ws.Cells.Item(132, 9).Value = 1270.3334  # I132: was 1270.6666
$ws.Cells.Item(132, 11).Value = 3811.0002  # K132: was 3811.9998
$ws.Cells.Item(132, 13).Value = -1281.0002  # M132: was -1281.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 3000  # H94: was 2904.5
$ws.Cells.Item(94, 9).Value = 3000  # I94: was 2904.5
$ws.Cells.Item(94, 11).Value = 3000  # K94: was 2904.5
$ws.Cells.Item(94, 13).Value = -2549  # M94: was -2453.5
# Row 97
$ws.Cells.Item(97, 8).Value = 9976.5  # H97: was 9981
$ws.Cells.Item(97, 9).Value = 9976.5  # I97: was 9981
$ws.Cells.Item(97, 11).Value = 9976.5  # K97: was 9981
$ws.Cells.Item(97, 13).Value = -8985.5  # M97: was -8990
# Row 99
$ws.Cells.Item(99, 8).Value = 1249.5  # H99: was 1249.75
$ws.Cells.Item(99, 10).Value = 1199  # J99: was 1199.5
$ws.Cells.Item(99, 12).Value = 1199  # L99: was 1199.5
$ws.Cells.Item(99, 14).Value = -4195  # N99: was -4195.5
# Row 107
$ws.Cells.Item(107, 8).Value = 1579.7  # H107: was 1589.2222
$ws.Cells.Item(107, 10).Value = 1497.4  # J107: was 1498.25
$ws.Cells.Item(107, 12).Value = 1497.4  # L107: was 1498.25
$ws.Cells.Item(107, 14).Value = -5337.4  # N107: was -5338.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 227.53847  # H7: was 245.75
$ws.Cells.Item(7, 9).Value = 79.77778000000001  # I7: was 88.625
$ws.Cells.Item(7, 11).Value = 79.77778000000001  # K7: was 88.625
$ws.Cells.Item(7, 13).Value = 33.22221999999999  # M7: was 24.375
# Row 99
$ws.Cells.Item(99, 8).Value = 4957.5  # H99: was 2183.3333
$ws.Cells.Item(99, 9).Value = 1600  # I99: was 0
$ws.Cells.Item(99, 10).Value = 6636.25  # J99: was 2183.3333
$ws.Cells.Item(99, 11).Value = 1600  # K99: was 0
$ws.Cells.Item(99, 12).Value = 6636.25  # L99: was 2183.3333
$ws.Cells.Item(99, 13).Value = -102  # M99: was None
$ws.Cells.Item(99, 14).Value = -9632.25  # N99: was -5179.3333
# Row 126
$ws.Cells.Item(126, 8).Value = 4957.5  # H126: was 2183.3333
$ws.Cells.Item(126, 9).Value = 1600  # I126: was 0
$ws.Cells.Item(126, 10).Value = 6636.25  # J126: was 2183.3333
$ws.Cells.Item(126, 11).Value = 4800  # K126: was 0
$ws.Cells.Item(126, 12).Value = 19908.75  # L126: was 6549.999899999999
$ws.Cells.Item(126, 13).Value = -2330  # M126: was None
$ws.Cells.Item(126, 14).Value = -24848.75  # N126: was -11489.9999
# Row 132
$ws.Cells.Item(132, 8).Value = 4233.625  # H132: was 4662
$ws.Cells.Item(132, 9).Value = 4267.4287  # I132: was 4795
$ws.Cells.Item(132, 11).Value = 12802.2861  # K132: was 14385
$ws.Cells.Item(132, 13).Value = -10272.2861  # M132: was -11855

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 132
$ws.Cells.Item(132, 8).Value = 1850  # H132: was 927.25
$ws.Cells.Item(132, 9).Value = 1850  # I132: was 1234.6666
$ws.Cells.Item(132, 10).Value = 0  # J132: was 5
$ws.Cells.Item(132, 11).Value = 16650  # K132: was 11111.9994
$ws.Cells.Item(132, 12).Value = 0  # L132: was 45
$ws.Cells.Item(132, 13).Value = -14120  # M132: was -8581.999400000001
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -5105

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 1972.5  # H80: was 1696.625
$ws.Cells.Item(80, 9).Value = 2057.5  # I80: was 1661.3334
$ws.Cells.Item(80, 11).Value = 2057.5  # K80: was 1661.3334
$ws.Cells.Item(80, 13).Value = -1059.5  # M80: was -663.3334
# Row 83
$ws.Cells.Item(83, 8).Value = 1972.5  # H83: was 1696.625
$ws.Cells.Item(83, 9).Value = 2057.5  # I83: was 1661.3334
$ws.Cells.Item(83, 11).Value = 10287.5  # K83: was 8306.666999999999
$ws.Cells.Item(83, 13).Value = -5295.5  # M83: was -3314.666999999999
# Row 126
$ws.Cells.Item(126, 8).Value = 6638  # H126: was 8957
$ws.Cells.Item(126, 9).Value = 2000  # I126: was 0
$ws.Cells.Item(126, 11).Value = 6000  # K126: was 0
$ws.Cells.Item(126, 13).Value = -3530  # M126: was None

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 3430.5  # H22: was 5200
$ws.Cells.Item(22, 9).Value = 3587.875  # I22: was 5200
$ws.Cells.Item(22, 10).Value = 2801  # J22: was 0
$ws.Cells.Item(22, 11).Value = 3587.875  # K22: was 5200
$ws.Cells.Item(22, 12).Value = 2801  # L22: was 0
$ws.Cells.Item(22, 13).Value = -3292.875  # M22: was -4905
$ws.Cells.Item(22, 14).Value = -3391  # N22: was None
# Row 27
$ws.Cells.Item(27, 8).Value = 3430.5  # H27: was 5200
$ws.Cells.Item(27, 9).Value = 3587.875  # I27: was 5200
$ws.Cells.Item(27, 10).Value = 2801  # J27: was 0
$ws.Cells.Item(27, 11).Value = 3587.875  # K27: was 5200
$ws.Cells.Item(27, 12).Value = 2801  # L27: was 0
$ws.Cells.Item(27, 13).Value = -3480.875  # M27: was -5093
$ws.Cells.Item(27, 14).Value = -3015  # N27: was None
# Row 46
$ws.Cells.Item(46, 8).Value = 1492.2972  # H46: was 1506.3715
$ws.Cells.Item(46, 9).Value = 1151.4736  # I46: was 1143.9
$ws.Cells.Item(46, 10).Value = 1852.0555  # J46: was 1989.6666
$ws.Cells.Item(46, 11).Value = 1151.4736  # K46: was 1143.9
$ws.Cells.Item(46, 12).Value = 1852.0555  # L46: was 1989.6666
$ws.Cells.Item(46, 13).Value = -963.4736  # M46: was -955.9000000000001
$ws.Cells.Item(46, 14).Value = -2228.0555  # N46: was -2365.6666
# Row 55
$ws.Cells.Item(55, 8).Value = 899.8261  # H55: was 872.75
$ws.Cells.Item(55, 10).Value = 1304.7778  # J55: was 1199.3
$ws.Cells.Item(55, 12).Value = 1304.7778  # L55: was 1199.3
$ws.Cells.Item(55, 14).Value = -1650.7778  # N55: was -1545.3
